$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.642.57"
$ws.Range("E2").Value = "  +1.11%  "

# Row 3
$ws.Range("D3").Value = "1.827.24"
$ws.Range("E3").Value = "  +1.92%  "

# Row 4
$ws.Range("E4").Value = "  +0.39%  "

# Row 5
$ws.Range("B5").Value = "USDC"
$ws.Range("C5").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.008"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.31%  "

# Row 6
$ws.Range("B6").Value = "BNB"
$ws.Range("C6").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "308.93"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.65%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4688"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.67%  "

# Row 8
$ws.Range("E8").Value = "  +0.18%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07142"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.00%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9022"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.06%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07670"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.19%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "19.41"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.12%  "

# Row 13
$ws.Range("D13").Value = "1.824.14"
$ws.Range("E13").Value = "  +1.31%  "

# Row 14
$ws.Range("E14").Value = "  -0.42%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.359"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.65%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "87.68"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.20%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.010"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.39%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008559"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.69%  "

# Row 19
$ws.Range("E19").Value = "  +0.39%  "

# Row 20
$ws.Range("D20").Value = "26.644.75"
$ws.Range("E20").Value = "  +1.01%  "

# Row 21
$ws.Range("E21").Value = "  -0.07%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.019"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.12%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.54"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.03%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.907"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.98%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.83"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.04%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "17.90"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.59%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.996"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.93%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "113.71"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.68%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.875"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.77%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08817"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.49%  "

# Row 31
$ws.Range("E31").Value = "  +1.19%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.846"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.70%  "

# Row 33
$ws.Range("E33").Value = "  +5.41%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7346"
$ws.Range("D34").Style = "Normal"

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.427"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.20%  "

# Row 36
$ws.Range("E36").Value = "  +0.84%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01926"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.10%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.951"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.34%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05145"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.44%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.856"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.44%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5051"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.13%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1498"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.00%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.066"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.70%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.009"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.32%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4647"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.58%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.04"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.52%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "98.35"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.27%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.570"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.46%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06022"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.04%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "63.69"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.05%  "

# Row 51
$ws.Range("E51").Value = "  -0.65%  "
